# Rename the three header/footer logo images:
#   - Pearson Edexcel logo (PearsonLogo.png) in both the "first page" footer
#     and the "default" footer: image1.png -> image2.png
#   - BTEC logo (BTec_Logo-Orange) in the "first page" header: image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers ---------------------------------------------------------
# wdHeaderFooterPrimary (1) and wdHeaderFooterFirstPage (2) both carry the
# Pearson logo picture; rename it from image1.png to image2.png.
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

# --- Headers -----------------------------------------------------------
# The "first page" header carries the BTEC logo; rename it from
# image2.jpg to image1.jpg.
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
